$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# The Huesca and Huelva rows trade places (row 53 and row 54):
#   Before: row 53 = Huelva (37 / 72 / 37 / 0)   row 54 = Huesca (37 / 0 / 37 / 0)
#   After:  row 53 = Huesca (37 / 0 / 37 / 0)    row 54 = Huelva (37 / 72 / 37 / 0)
# Columns B, D and E are identical between the two rows, so only the city
# name (column A) and "Casos activos" (column C) actually need updating.
$ws.Range("A53").Value = "Huesca"
$ws.Range("C53").Value = 0

$ws.Range("A54").Value = "Huelva"
$ws.Range("C54").Value = 72

# Refresh the "last updated" timestamp shown in A1.
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 21:16"
